$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 596678
$ws.Range("J43").Value = 1174123.8
$ws.Range("L43").Value = 1174123.8
$ws.Range("N43").Value = -1174261.8
$ws.Range("H98").Value = 1975.8837
$ws.Range("I98").Value = 2112.6667
$ws.Range("K98").Value = 2112.6667
$ws.Range("M98").Value = -614.6667000000002
$ws.Range("H111").Value = 10418250
$ws.Range("I111").Value = 13890366
$ws.Range("K111").Value = 41671098
$ws.Range("M111").Value = -41668031
$ws.Range("H112").Value = 9430.825999999999
$ws.Range("J112").Value = 10107.333
$ws.Range("L112").Value = 30321.999
$ws.Range("N112").Value = -32537.999
$ws.Range("H122").Value = 1975.8837
$ws.Range("I122").Value = 2112.6667
$ws.Range("K122").Value = 6338.000100000001
$ws.Range("M122").Value = -3888.000100000001
$ws.Range("H135").Value = 182432.34
$ws.Range("I135").Value = 204661.55
$ws.Range("K135").Value = 1841953.95
$ws.Range("M135").Value = -1839418.95
$ws.Range("H137").Value = 2678.875
$ws.Range("I137").Value = 1750
$ws.Range("J137").Value = 2988.5
$ws.Range("K137").Value = 5250
$ws.Range("L137").Value = 8965.5
$ws.Range("M137").Value = -2700
$ws.Range("N137").Value = -14065.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2534.3333
$ws.Range("I2").Value = 2189
$ws.Range("J2").Value = 2929
$ws.Range("K2").Value = 2189
$ws.Range("L2").Value = 2929
$ws.Range("M2").Value = -2076
$ws.Range("N2").Value = -3155
$ws.Range("H32").Value = 2663584.8
$ws.Range("I32").Value = 2910058.2
$ws.Range("K32").Value = 2910058.2
$ws.Range("M32").Value = -2909771.2
$ws.Range("H45").Value = 2647.2856
$ws.Range("I45").Value = 1780.25
$ws.Range("K45").Value = 1780.25
$ws.Range("M45").Value = -1403.25
$ws.Range("H74").Value = 71991.734
$ws.Range("I74").Value = 114652.89
$ws.Range("K74").Value = 114652.89
$ws.Range("M74").Value = -113778.89
$ws.Range("H77").Value = 71991.734
$ws.Range("I77").Value = 114652.89
$ws.Range("K77").Value = 573264.45
$ws.Range("M77").Value = -568896.45
$ws.Range("H116").Value = 2534.3333
$ws.Range("I116").Value = 2189
$ws.Range("J116").Value = 2929
$ws.Range("K116").Value = 2189
$ws.Range("L116").Value = 2929
$ws.Range("M116").Value = 105
$ws.Range("N116").Value = -7517
$ws.Range("H132").Value = 1674617.4
$ws.Range("I132").Value = 3854091.5
$ws.Range("K132").Value = 11562274.5
$ws.Range("M132").Value = -11559744.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2534.3333
$ws.Range("I3").Value = 2189
$ws.Range("J3").Value = 2929
$ws.Range("K3").Value = 2189
$ws.Range("L3").Value = 2929
$ws.Range("M3").Value = -2075
$ws.Range("N3").Value = -3157
$ws.Range("H99").Value = 6063921.5
$ws.Range("I99").Value = 2811
$ws.Range("K99").Value = 2811
$ws.Range("M99").Value = -1313
$ws.Range("H105").Value = 3986.111
$ws.Range("I105").Value = 2760
$ws.Range("J105").Value = 5212.222
$ws.Range("K105").Value = 2760
$ws.Range("L105").Value = 5212.222
$ws.Range("M105").Value = -1013
$ws.Range("N105").Value = -8706.222
$ws.Range("H134").Value = 7687.92
$ws.Range("I134").Value = 2974
$ws.Range("J134").Value = 10339.5
$ws.Range("K134").Value = 8922
$ws.Range("L134").Value = 31018.5
$ws.Range("M134").Value = -6387
$ws.Range("N134").Value = -36088.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3676.9375
$ws.Range("I16").Value = 3229.375
$ws.Range("K16").Value = 3229.375
$ws.Range("M16").Value = -2942.375
$ws.Range("H31").Value = 8429.423000000001
$ws.Range("I31").Value = 1210.375
$ws.Range("J31").Value = 11637.889
$ws.Range("K31").Value = 1210.375
$ws.Range("L31").Value = 11637.889
$ws.Range("M31").Value = -915.375
$ws.Range("N31").Value = -12227.889
$ws.Range("H34").Value = 8429.423000000001
$ws.Range("I34").Value = 1210.375
$ws.Range("J34").Value = 11637.889
$ws.Range("K34").Value = 1210.375
$ws.Range("L34").Value = 11637.889
$ws.Range("M34").Value = -1008.375
$ws.Range("N34").Value = -12041.889
$ws.Range("H107").Value = 1352.4166
$ws.Range("I107").Value = 620.8461
$ws.Range("J107").Value = 2217
$ws.Range("K107").Value = 620.8461
$ws.Range("L107").Value = 2217
$ws.Range("M107").Value = 1299.1539
$ws.Range("N107").Value = -6057
$ws.Range("H109").Value = 53127.668
$ws.Range("J109").Value = 53127.668
$ws.Range("L109").Value = 53127.668
$ws.Range("N109").Value = -55207.668
$ws.Range("H113").Value = 3676.9375
$ws.Range("I113").Value = 3229.375
$ws.Range("K113").Value = 3229.375
$ws.Range("M113").Value = -1059.375

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2148.2334
$ws.Range("I68").Value = 1493.8889
$ws.Range("J68").Value = 2428.6667
$ws.Range("K68").Value = 4481.6667
$ws.Range("L68").Value = 7286.000100000001
$ws.Range("M68").Value = -3670.6667
$ws.Range("N68").Value = -8908.000100000001
$ws.Range("H71").Value = 2148.2334
$ws.Range("I71").Value = 1493.8889
$ws.Range("J71").Value = 2428.6667
$ws.Range("K71").Value = 13445.0001
$ws.Range("L71").Value = 21858.0003
$ws.Range("M71").Value = -9389.000099999999
$ws.Range("N71").Value = -29970.0003
$ws.Range("H108").Value = 1263.3334
$ws.Range("I108").Value = 995
$ws.Range("J108").Value = 1800
$ws.Range("K108").Value = 2985
$ws.Range("L108").Value = 5400
$ws.Range("M108").Value = -105
$ws.Range("N108").Value = -11160
$ws.Range("H131").Value = 1969.36
$ws.Range("J131").Value = 2159.6287
$ws.Range("L131").Value = 6478.886100000001
$ws.Range("N131").Value = -16558.8861

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 169666.33
$ws.Range("J80").Value = 169666.33
$ws.Range("L80").Value = 169666.33
$ws.Range("N80").Value = -171662.33
$ws.Range("H83").Value = 169666.33
$ws.Range("J83").Value = 169666.33
$ws.Range("L83").Value = 848331.6499999999
$ws.Range("N83").Value = -858315.6499999999
$ws.Range("H113").Value = 6535.6763
$ws.Range("I113").Value = 3900
$ws.Range("K113").Value = 3900
$ws.Range("M113").Value = -1730
$ws.Range("H122").Value = 143002820
$ws.Range("I122").Value = 166835790
$ws.Range("J122").Value = 5008
$ws.Range("K122").Value = 500507370
$ws.Range("L122").Value = 15024
$ws.Range("M122").Value = -500504920
$ws.Range("N122").Value = -19924
$ws.Range("H126").Value = 7199.857
$ws.Range("J126").Value = 7199.857
$ws.Range("L126").Value = 21599.571
$ws.Range("N126").Value = -26539.571
$ws.Range("H132").Value = 6580.846
$ws.Range("J132").Value = 7950.1113
$ws.Range("L132").Value = 23850.3339
$ws.Range("N132").Value = -28910.3339
$ws.Range("H134").Value = 91698.164
$ws.Range("J134").Value = 91698.164
$ws.Range("L134").Value = 275094.492
$ws.Range("N134").Value = -280164.492

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7864.091
$ws.Range("I7").Value = 7374.5
$ws.Range("J7").Value = 8143.857
$ws.Range("K7").Value = 7374.5
$ws.Range("L7").Value = 8143.857
$ws.Range("M7").Value = -7262.5
$ws.Range("N7").Value = -8367.857
$ws.Range("H55").Value = 1024
$ws.Range("I55").Value = 1168.6154
$ws.Range("K55").Value = 1168.6154
$ws.Range("M55").Value = -995.6153999999999
$ws.Range("H87").Value = 60000
$ws.Range("J87").Value = 60000
$ws.Range("L87").Value = 60000
$ws.Range("N87").Value = -62246
$ws.Range("H90").Value = 60000
$ws.Range("J90").Value = 60000
$ws.Range("L90").Value = 180000
$ws.Range("N90").Value = -191232
$ws.Range("H122").Value = 8484.409
$ws.Range("I122").Value = 9077.532999999999
$ws.Range("K122").Value = 27232.599
$ws.Range("M122").Value = -24782.599
$ws.Range("H126").Value = 7864.091
$ws.Range("I126").Value = 7374.5
$ws.Range("J126").Value = 8143.857
$ws.Range("K126").Value = 22123.5
$ws.Range("L126").Value = 24431.571
$ws.Range("M126").Value = -19653.5
$ws.Range("N126").Value = -29371.571
$ws.Range("H136").Value = 8472.200000000001
$ws.Range("I136").Value = 5403.625
$ws.Range("J136").Value = 11979.143
$ws.Range("K136").Value = 16210.875
$ws.Range("L136").Value = 35937.429
$ws.Range("M136").Value = -13660.875
$ws.Range("N136").Value = -41037.429

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 893.6667
$ws.Range("I107").Value = 1067.3334
$ws.Range("K107").Value = 3202.0002
$ws.Range("M107").Value = -1282.0002
$ws.Range("H113").Value = 32810
$ws.Range("I113").Value = 63287.25
$ws.Range("K113").Value = 189861.75
$ws.Range("M113").Value = -187691.75
$ws.Range("H122").Value = 147164.89
$ws.Range("I122").Value = 177879
$ws.Range("J122").Value = 5880
$ws.Range("K122").Value = 533637
$ws.Range("L122").Value = 17640
$ws.Range("M122").Value = -531187
$ws.Range("N122").Value = -22540
$ws.Range("H132").Value = 13904.225
$ws.Range("I132").Value = 9977.833000000001
$ws.Range("K132").Value = 29933.499
$ws.Range("M132").Value = -27403.499
$ws.Range("H136").Value = 29097.95
$ws.Range("I136").Value = 1026.92
$ws.Range("K136").Value = 3080.76
$ws.Range("M136").Value = -530.7600000000002
